$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 12560.3
$ws.Range("I32").Value = 19533
$ws.Range("J32").Value = 9572
$ws.Range("K32").Value = 19533
$ws.Range("L32").Value = 9572
$ws.Range("M32").Value = -19207
$ws.Range("N32").Value = -10224
$ws.Range("H80").Value = 1287.1666
$ws.Range("I80").Value = 372.9091
$ws.Range("K80").Value = 1118.7273
$ws.Range("M80").Value = -120.7273
$ws.Range("H83").Value = 1287.1666
$ws.Range("I83").Value = 372.9091
$ws.Range("K83").Value = 3356.1819
$ws.Range("M83").Value = 1635.8181
$ws.Range("H111").Value = 2389.9473
$ws.Range("I111").Value = 2320.4375
$ws.Range("K111").Value = 6961.3125
$ws.Range("M111").Value = -3894.3125
$ws.Range("H112").Value = 2818.0667
$ws.Range("J112").Value = 2876.5
$ws.Range("L112").Value = 8629.5
$ws.Range("N112").Value = -10845.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H132").Value = 1040.1875
$ws.Range("I132").Value = 663.6070999999999
$ws.Range("K132").Value = 1990.8213
$ws.Range("M132").Value = 539.1787000000002
$ws.Range("H137").Value = 11234.863
$ws.Range("I137").Value = 5054.048
$ws.Range("J137").Value = 16878.217
$ws.Range("K137").Value = 15162.144
$ws.Range("L137").Value = 50634.651
$ws.Range("M137").Value = -12612.144
$ws.Range("N137").Value = -55734.651
$ws.Range("H138").Value = 3516.7659
$ws.Range("I138").Value = 3272.2
$ws.Range("J138").Value = 4914.2856
$ws.Range("K138").Value = 9816.599999999999
$ws.Range("L138").Value = 14742.8568
$ws.Range("M138").Value = -4676.599999999999
$ws.Range("N138").Value = -25022.8568
$ws.Range("H141").Value = 1500.1538
$ws.Range("I141").Value = 1380.16
$ws.Range("K141").Value = 4140.48
$ws.Range("M141").Value = 1039.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15171386
$ws.Range("I32").Value = 16413537
$ws.Range("J32").Value = 17135.6
$ws.Range("K32").Value = 16413537
$ws.Range("L32").Value = 17135.6
$ws.Range("M32").Value = -16413250
$ws.Range("N32").Value = -17709.6
$ws.Range("H74").Value = 2766.4814
$ws.Range("I74").Value = 1931.409
$ws.Range("K74").Value = 1931.409
$ws.Range("M74").Value = -1057.409
$ws.Range("H77").Value = 2766.4814
$ws.Range("I77").Value = 1931.409
$ws.Range("K77").Value = 9657.045
$ws.Range("M77").Value = -5289.045
$ws.Range("H122").Value = 1788.0476
$ws.Range("I122").Value = 1596.25
$ws.Range("K122").Value = 4788.75
$ws.Range("M122").Value = -2338.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2844.9565
$ws.Range("J20").Value = 2973.6365
$ws.Range("L20").Value = 2973.6365
$ws.Range("N20").Value = -3467.6365
$ws.Range("H105").Value = 20007084
$ws.Range("I105").Value = 31259154
$ws.Range("K105").Value = 31259154
$ws.Range("M105").Value = -31257407
$ws.Range("H107").Value = 926.25714
$ws.Range("I107").Value = 837.9259
$ws.Range("J107").Value = 1224.375
$ws.Range("K107").Value = 837.9259
$ws.Range("L107").Value = 1224.375
$ws.Range("M107").Value = 1082.0741
$ws.Range("N107").Value = -5064.375
$ws.Range("H124").Value = 66333
$ws.Range("J124").Value = 66333
$ws.Range("L124").Value = 66333
$ws.Range("N124").Value = -76153
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3318.5
$ws.Range("I16").Value = 1727.75
$ws.Range("K16").Value = 1727.75
$ws.Range("M16").Value = -1440.75
$ws.Range("H28").Value = 25000
$ws.Range("J28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("N28").Value = -25490
$ws.Range("H31").Value = 4046.7454
$ws.Range("I31").Value = 1832.591
$ws.Range("K31").Value = 1832.591
$ws.Range("M31").Value = -1537.591
$ws.Range("H34").Value = 4046.7454
$ws.Range("I34").Value = 1832.591
$ws.Range("K34").Value = 1832.591
$ws.Range("M34").Value = -1630.591
$ws.Range("H58").Value = 3064.98
$ws.Range("I58").Value = 2405.125
$ws.Range("J58").Value = 4238.0557
$ws.Range("K58").Value = 2405.125
$ws.Range("L58").Value = 4238.0557
$ws.Range("M58").Value = -2202.125
$ws.Range("N58").Value = -4644.0557
$ws.Range("H99").Value = 4936
$ws.Range("J99").Value = 2996.6667
$ws.Range("L99").Value = 2996.6667
$ws.Range("N99").Value = -5992.6667
$ws.Range("H113").Value = 3318.5
$ws.Range("I113").Value = 1727.75
$ws.Range("K113").Value = 1727.75
$ws.Range("M113").Value = 442.25
$ws.Range("H123").Value = 59999.5
$ws.Range("J123").Value = 59999.5
$ws.Range("L123").Value = 59999.5
$ws.Range("N123").Value = -69799.5
$ws.Range("H126").Value = 4936
$ws.Range("J126").Value = 2996.6667
$ws.Range("L126").Value = 8990.000100000001
$ws.Range("N126").Value = -13930.0001
$ws.Range("H127").Value = 60000
$ws.Range("J127").Value = 60000
$ws.Range("L127").Value = 60000
$ws.Range("N127").Value = -69920
$ws.Range("H132").Value = 2431.5366
$ws.Range("I132").Value = 1775.7812
$ws.Range("J132").Value = 4763.1113
$ws.Range("K132").Value = 5327.3436
$ws.Range("L132").Value = 14289.3339
$ws.Range("M132").Value = -2797.3436
$ws.Range("N132").Value = -19349.3339
$ws.Range("H134").Value = 3606.973
$ws.Range("I134").Value = 2981.6858
$ws.Range("K134").Value = 8945.057400000002
$ws.Range("M134").Value = -6410.057400000002
$ws.Range("H136").Value = 3064.98
$ws.Range("I136").Value = 2405.125
$ws.Range("J136").Value = 4238.0557
$ws.Range("K136").Value = 7215.375
$ws.Range("L136").Value = 12714.1671
$ws.Range("M136").Value = -4665.375
$ws.Range("N136").Value = -17814.1671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2375
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H131").Value = 9356.706
$ws.Range("J131").Value = 9819
$ws.Range("L131").Value = 29457
$ws.Range("N131").Value = -39537
$ws.Range("I137").Value = 949.625
$ws.Range("J137").Value = 2105
$ws.Range("K137").Value = 2848.875
$ws.Range("L137").Value = 6315
$ws.Range("M137").Value = 2251.125
$ws.Range("N137").Value = -16515

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4142.9546
$ws.Range("J70").Value = 3956.75
$ws.Range("L70").Value = 3956.75
$ws.Range("N70").Value = -4496.75
$ws.Range("H73").Value = 4142.9546
$ws.Range("J73").Value = 3956.75
$ws.Range("L73").Value = 3956.75
$ws.Range("N73").Value = -5828.75
$ws.Range("H122").Value = 1805.0435
$ws.Range("J122").Value = 2214.2
$ws.Range("L122").Value = 6642.599999999999
$ws.Range("N122").Value = -11542.6
$ws.Range("H132").Value = 4191.0625
$ws.Range("I132").Value = 3004.0715
$ws.Range("K132").Value = 9012.2145
$ws.Range("M132").Value = -6482.2145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 45000
$ws.Range("J69").Value = 45000
$ws.Range("L69").Value = 45000
$ws.Range("N69").Value = -46622
$ws.Range("H72").Value = 45000
$ws.Range("J72").Value = 45000
$ws.Range("L72").Value = 135000
$ws.Range("N72").Value = -143112
$ws.Range("H136").Value = 3684.0645
$ws.Range("I136").Value = 3252.5217
$ws.Range("K136").Value = 9757.5651
$ws.Range("M136").Value = -7207.5651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 64875
$ws.Range("I17").Value = 17250
$ws.Range("J17").Value = 112500
$ws.Range("K17").Value = 17250
$ws.Range("L17").Value = 112500
$ws.Range("M17").Value = -17078
$ws.Range("N17").Value = -112844
$ws.Range("H100").Value = 598
$ws.Range("I100").Value = 419.2
$ws.Range("K100").Value = 838.4
$ws.Range("M100").Value = -297.4
$ws.Range("H126").Value = 3688.6333
$ws.Range("I126").Value = 3773.95
$ws.Range("K126").Value = 11321.85
$ws.Range("M126").Value = -8851.849999999999
$ws.Range("H131").Value = 119664.29
$ws.Range("J131").Value = 126166.664
$ws.Range("L131").Value = 126166.664
$ws.Range("N131").Value = -136246.664
$ws.Range("H132").Value = 4756.1763
$ws.Range("I132").Value = 3487.1365
$ws.Range("K132").Value = 10461.4095
$ws.Range("M132").Value = -7931.4095
$ws.Range("H133").Value = 108637.5
$ws.Range("J133").Value = 108637.5
$ws.Range("L133").Value = 108637.5
$ws.Range("N133").Value = -118757.5
$ws.Range("H136").Value = 1685.1333
$ws.Range("I136").Value = 569.1799999999999
$ws.Range("K136").Value = 1707.54
$ws.Range("M136").Value = 842.46
